# Update the division-problem worksheet numbers (table of three-digit / one-digit
# division problems). Each Find/Replace targets the exact, unique "old" expression
# text in the single run of the matching table cell and swaps in the new problem.
# wdReplaceOne (1) is used throughout, and replacements are applied in document
# order, so the transient reuse of "423÷4=" as a later target value never collides
# with an earlier, not-yet-updated cell.
$d = $word.ActiveDocument

$d.Content.Find.Execute("696÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "635÷2=", 1) | Out-Null
$d.Content.Find.Execute("154÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "239÷4=", 1) | Out-Null
$d.Content.Find.Execute("423÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "911÷8=", 1) | Out-Null
$d.Content.Find.Execute("212÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "944÷7=", 1) | Out-Null
$d.Content.Find.Execute("835÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "948÷4=", 1) | Out-Null
$d.Content.Find.Execute("699÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "199÷4=", 1) | Out-Null
$d.Content.Find.Execute("293÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "693÷7=", 1) | Out-Null
$d.Content.Find.Execute("497÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "680÷9=", 1) | Out-Null
$d.Content.Find.Execute("269÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "296÷2=", 1) | Out-Null
$d.Content.Find.Execute("354÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "401÷2=", 1) | Out-Null
$d.Content.Find.Execute("331÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "569÷4=", 1) | Out-Null
$d.Content.Find.Execute("539÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "913÷4=", 1) | Out-Null
$d.Content.Find.Execute("187÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "820÷2=", 1) | Out-Null
$d.Content.Find.Execute("260÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "756÷5=", 1) | Out-Null
$d.Content.Find.Execute("780÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "339÷9=", 1) | Out-Null
$d.Content.Find.Execute("321÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "297÷2=", 1) | Out-Null
$d.Content.Find.Execute("605÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "177÷8=", 1) | Out-Null
$d.Content.Find.Execute("114÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "273÷9=", 1) | Out-Null
$d.Content.Find.Execute("586÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "289÷3=", 1) | Out-Null
$d.Content.Find.Execute("498÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "921÷8=", 1) | Out-Null
$d.Content.Find.Execute("240÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "496÷7=", 1) | Out-Null
$d.Content.Find.Execute("296÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "836÷8=", 1) | Out-Null
$d.Content.Find.Execute("418÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "694÷7=", 1) | Out-Null
$d.Content.Find.Execute("978÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "423÷4=", 1) | Out-Null
$d.Content.Find.Execute("646÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "342÷7=", 1) | Out-Null

Write-Output "Updated 25 division problems."
